# Auto-generated COM-interop script applying the Tonberry_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H33").Value = 107
$ws.Range("I33").Value = 84.5
$ws.Range("K33").Value = 84.5
$ws.Range("M33").Value = 144.5
$ws.Range("H106").Value = 3693.1538
$ws.Range("I106").Value = 2444.3333
$ws.Range("J106").Value = 6503
$ws.Range("K106").Value = 2444.3333
$ws.Range("L106").Value = 6503
$ws.Range("M106").Value = -1813.3333
$ws.Range("N106").Value = -7765
$ws.Range("H138").Value = 2557
$ws.Range("I138").Value = 2870.5454
$ws.Range("K138").Value = 8611.636200000001
$ws.Range("M138").Value = -3471.636200000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4002.2104
$ws.Range("I32").Value = 2604.634
$ws.Range("K32").Value = 2604.634
$ws.Range("M32").Value = -2317.634
$ws.Range("H45").Value = 1672.7778
$ws.Range("I45").Value = 1233
$ws.Range("K45").Value = 1233
$ws.Range("M45").Value = -856
$ws.Range("H61").Value = 2861.05
$ws.Range("I61").Value = 2329.4666
$ws.Range("J61").Value = 4455.8
$ws.Range("K61").Value = 2329.4666
$ws.Range("L61").Value = 4455.8
$ws.Range("M61").Value = -2117.4666
$ws.Range("N61").Value = -4879.8
$ws.Range("H122").Value = 1717.8667
$ws.Range("I122").Value = 1478.2858
$ws.Range("K122").Value = 4434.857400000001
$ws.Range("M122").Value = -1984.857400000001
$ws.Range("H132").Value = 2617.1428
$ws.Range("I132").Value = 1455.5
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 4366.5
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = -1836.5
$ws.Range("N132").Value = -17558
$ws.Range("H136").Value = 2861.05
$ws.Range("I136").Value = 2329.4666
$ws.Range("J136").Value = 4455.8
$ws.Range("K136").Value = 6988.399800000001
$ws.Range("L136").Value = 13367.4
$ws.Range("M136").Value = -4438.399800000001
$ws.Range("N136").Value = -18467.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 135261.67
$ws.Range("I86").Value = 1796.5
$ws.Range("J86").Value = 183794.45
$ws.Range("K86").Value = 1796.5
$ws.Range("L86").Value = 183794.45
$ws.Range("M86").Value = -673.5
$ws.Range("N86").Value = -186040.45
$ws.Range("H89").Value = 135261.67
$ws.Range("I89").Value = 1796.5
$ws.Range("J89").Value = 183794.45
$ws.Range("K89").Value = 8982.5
$ws.Range("L89").Value = 918972.25
$ws.Range("M89").Value = -3366.5
$ws.Range("N89").Value = -930204.25
$ws.Range("H99").Value = 996
$ws.Range("I99").Value = 996
$ws.Range("K99").Value = 996
$ws.Range("M99").Value = 502
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 764.7143
$ws.Range("I16").Value = 727.3333
$ws.Range("K16").Value = 727.3333
$ws.Range("M16").Value = -440.3333
$ws.Range("H58").Value = 2289613.8
$ws.Range("I58").Value = 3345567
$ws.Range("J58").Value = 1715.3334
$ws.Range("K58").Value = 3345567
$ws.Range("L58").Value = 1715.3334
$ws.Range("M58").Value = -3345364
$ws.Range("N58").Value = -2121.3334
$ws.Range("H94").Value = 893.5714
$ws.Range("I94").Value = 764.8333
$ws.Range("J94").Value = 990.125
$ws.Range("K94").Value = 764.8333
$ws.Range("L94").Value = 990.125
$ws.Range("M94").Value = -313.8333
$ws.Range("N94").Value = -1892.125
$ws.Range("H99").Value = 2321.4285
$ws.Range("I99").Value = 1966.6666
$ws.Range("J99").Value = 2587.5
$ws.Range("K99").Value = 1966.6666
$ws.Range("L99").Value = 2587.5
$ws.Range("M99").Value = -468.6666
$ws.Range("N99").Value = -5583.5
$ws.Range("H113").Value = 764.7143
$ws.Range("I113").Value = 727.3333
$ws.Range("K113").Value = 727.3333
$ws.Range("M113").Value = 1442.6667
$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -49910
$ws.Range("H126").Value = 2321.4285
$ws.Range("I126").Value = 1966.6666
$ws.Range("J126").Value = 2587.5
$ws.Range("K126").Value = 5899.9998
$ws.Range("L126").Value = 7762.5
$ws.Range("M126").Value = -3429.9998
$ws.Range("N126").Value = -12702.5
$ws.Range("H134").Value = 854.0909
$ws.Range("I134").Value = 843.8889
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 2531.6667
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = 3.333299999999781
$ws.Range("N134").Value = -7770
$ws.Range("H136").Value = 2289613.8
$ws.Range("I136").Value = 3345567
$ws.Range("J136").Value = 1715.3334
$ws.Range("K136").Value = 10036701
$ws.Range("L136").Value = 5146.0002
$ws.Range("M136").Value = -10034151
$ws.Range("N136").Value = -10246.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 48.5
$ws.Range("I33").Value = 48.333332
$ws.Range("J33").Value = 49
$ws.Range("K33").Value = 289.999992
$ws.Range("L33").Value = 294
$ws.Range("M33").Value = -6.99999200000002
$ws.Range("N33").Value = -860
$ws.Range("H139").Value = 11240.909
$ws.Range("I139").Value = 12165.2
$ws.Range("K139").Value = 36495.60000000001
$ws.Range("M139").Value = -31355.60000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2695439.5
$ws.Range("I126").Value = 4275697
$ws.Range("J126").Value = 127521.75
$ws.Range("K126").Value = 12827091
$ws.Range("L126").Value = 382565.25
$ws.Range("M126").Value = -12824621
$ws.Range("N126").Value = -387505.25
$ws.Range("H132").Value = 5499283.5
$ws.Range("I132").Value = 38461536
$ws.Range("K132").Value = 115384608
$ws.Range("M132").Value = -115382078
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H68").Value = 5159.8335
$ws.Range("I68").Value = 5494.75
$ws.Range("K68").Value = 5494.75
$ws.Range("M68").Value = -4745.75
$ws.Range("H71").Value = 5159.8335
$ws.Range("I71").Value = 5494.75
$ws.Range("K71").Value = 27473.75
$ws.Range("M71").Value = -23729.75
$ws.Range("H82").Value = 3997.25
$ws.Range("I82").Value = 1333
$ws.Range("J82").Value = 5595.8
$ws.Range("K82").Value = 1333
$ws.Range("L82").Value = 5595.8
$ws.Range("M82").Value = -972
$ws.Range("N82").Value = -6317.8
$ws.Range("H85").Value = 3997.25
$ws.Range("I85").Value = 1333
$ws.Range("J85").Value = 5595.8
$ws.Range("K85").Value = 1333
$ws.Range("L85").Value = 5595.8
$ws.Range("M85").Value = -85
$ws.Range("N85").Value = -8091.8
$ws.Range("H132").Value = 2142.6667
$ws.Range("I132").Value = 1399.4286
$ws.Range("J132").Value = 3183.2
$ws.Range("K132").Value = 4198.2858
$ws.Range("L132").Value = 9549.599999999999
$ws.Range("M132").Value = -1668.2858
$ws.Range("N132").Value = -14609.6
$ws.Range("H136").Value = 2749.861
$ws.Range("J136").Value = 4758.5
$ws.Range("L136").Value = 14275.5
$ws.Range("N136").Value = -19375.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 28932.334
$ws.Range("J92").Value = 28932.334
$ws.Range("L92").Value = 28932.334
$ws.Range("N92").Value = -33924.334
$ws.Range("H107").Value = 992.1667
$ws.Range("I107").Value = 537.5
$ws.Range("J107").Value = 1901.5
$ws.Range("K107").Value = 1612.5
$ws.Range("L107").Value = 5704.5
$ws.Range("M107").Value = 307.5
$ws.Range("N107").Value = -9544.5
$ws.Range("H132").Value = 5018.467
$ws.Range("I132").Value = 1062.35
$ws.Range("J132").Value = 8183.36
$ws.Range("K132").Value = 3187.05
$ws.Range("L132").Value = 24550.08
$ws.Range("M132").Value = -657.0499999999997
$ws.Range("N132").Value = -29610.08
$ws.Range("H136").Value = 16837204
$ws.Range("I136").Value = 26457044
$ws.Range("K136").Value = 79371132
$ws.Range("M136").Value = -79368582
